# Reorder the "Recorded By" (column G) list of names/emails so that the
# literal token "System" (exact case) is moved from the front of the
# comma-separated list to the end, preserving the relative order of the
# remaining entries (i.e. the whole list is reversed, since "System" was
# always the first entry when present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$colG = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -ceq "System") {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $count = $parts.Length
        $reversedParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversedParts)
        $cell.Value = $newText
    }
}
